# Generate Report for Handoff
# Updates the localization-status workbook so that the "dea62695..." file
# (row 3 in every sheet) is marked as "Ready for handoff" instead of
# "Handed back: in sync with en-US", refreshes its handoff timestamps and
# records an error detail describing that the handback file is stale.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dfda8bc0095100fac833f22a54768b87e61c9234/e2e/dea62695-3c19-4b4d-8586-a527fcf8cb49.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c0270a0f63e64e92dd7b5457cafb54dd1d7503c/e2e/dea62695-3c19-4b4d-8586-a527fcf8cb49.md."

# --- "Overview" sheet: zh-cn / de-de status + latest generate date for the
#     dea62695 file (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-01 12:54:26"

# --- "zh-cn" sheet: Status / Latest Handoff Datetime / Error Detail for the
#     dea62695 file (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-09-01 12:54:21"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- "de-de" sheet: Status / Latest Handoff Datetime / Error Detail for the
#     dea62695 file (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-09-01 12:54:26"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
